$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Set column D (rows 2-8) to "A" (row 9 already has "A")
$ws.Range("D2").Value = "A"
$ws.Range("D3").Value = "A"
$ws.Range("D4").Value = "A"
$ws.Range("D5").Value = "A"
$ws.Range("D6").Value = "A"
$ws.Range("D7").Value = "A"
$ws.Range("D8").Value = "A"

# Clear column J (rows 2-9), previously "Pass"/"Fail"
$ws.Range("J2").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("J9").Value = ""

# Clear the active selection marker left on the sheet view
$ws.Range("A1").Select()
